$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "user_payment_id" column (H) with header + sample value,
# mirroring the other missing-record columns already on the sheet.
$ws.Range("H1").Value = "user_payment_id"
$ws.Range("H2").Value = "UP_MX_hotgo_95345765_1605768224"

# H1 gets the same header formatting as the other header cells (G1).
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# H2 gets a centered, wrapped style (new cellXfs entry).
$ws.Range("H2").HorizontalAlignment = -4108
$ws.Range("H2").VerticalAlignment = -4108
$ws.Range("H2").WrapText = $true

# Widen the new column to fit the sample value.
$ws.Range("H1:H9").ColumnWidth = 36.109375

# Match the author's final selection.
$ws.Range("H13").Select() | Out-Null
